$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the flavor_text for super_contest_effect_id = 4 (row 4) to use the
# doubled apostrophe form of the text.
$ws.Range("C4").Value = "Earn +2 if the Judge''s Voltage goes up."
